# Rename the three header/footer logo pictures (their docPr/cNvPr "name"
# attribute only -- the pictures themselves, their ids and their
# descriptions are left untouched):
#
#   footer "first page"  (footer1.xml, docPr id=3): image2.png -> image1.png
#   footer "default"     (footer2.xml, docPr id=2): image2.png -> image1.png
#   header "first page"  (header1.xml, docPr id=1): image1.jpg -> image2.jpg
#
# InlineShape has no settable .Name in the Word object model, so each
# picture is briefly converted to a floating Shape (where .Name is
# writable), renamed, and converted straight back to an inline shape so
# the on-disk <wp:inline .../> wrapping is preserved.

$d = $word.ActiveDocument
$section = $d.Sections(1)

function Rename-InlineLogo($inlineShape, $newName) {
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

# wdHeaderFooterFirstPage = 2, wdHeaderFooterPrimary = 1
$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2

# --- footer1.xml : the "first page" footer's Pearson logo ---
$footerFirst = $section.Footers($wdHeaderFooterFirstPage)
$pearsonLogoFirstFooter = $footerFirst.Range.InlineShapes(1)
Rename-InlineLogo $pearsonLogoFirstFooter "image1.png"

# --- footer2.xml : the "default" (primary) footer's Pearson logo ---
$footerDefault = $section.Footers($wdHeaderFooterPrimary)
$pearsonLogoDefaultFooter = $footerDefault.Range.InlineShapes(1)
Rename-InlineLogo $pearsonLogoDefaultFooter "image1.png"

# --- header1.xml : the "first page" header's BTEC logo ---
$headerFirst = $section.Headers($wdHeaderFooterFirstPage)
$btecLogoFirstHeader = $headerFirst.Range.InlineShapes(1)
Rename-InlineLogo $btecLogoFirstHeader "image2.jpg"
